$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1 (0-based row 0): 99.89 -> 0M
$t.Cell(1,1).Range.Text = "0M"

# Row 2 (0-based row 1): 0.1 -> 0M
$t.Cell(2,1).Range.Text = "0M"

# Row 3 (0-based row 2): 90 -> 0M
$t.Cell(3,1).Range.Text = "0M"

# Row 4 (0-based row 3): 103 -> 412
$t.Cell(4,1).Range.Text = "412"

# Row 6 (0-based row 5): 0.00027 -> 0.00091
$t.Cell(6,1).Range.Text = "0.00091"

# Row 7 (0-based row 6): 0.00011 -> 0.00024
$t.Cell(7,1).Range.Text = "0.00024"

# Row 8 (0-based row 7): 0.00003 -> 0.00007
$t.Cell(8,1).Range.Text = "0.00007"

# Row 9 (0-based row 8): 0.00010 -> 0.00040
$t.Cell(9,1).Range.Text = "0.00040"

# Row 10 (0-based row 9): 0.00012 -> 0.00044
$t.Cell(10,1).Range.Text = "0.00044"

# Row 11 (0-based row 10): 0.00012 -> 0.00052
$t.Cell(11,1).Range.Text = "0.00052"

# Row 12 (0-based row 11): 0.01148 -> 0.09883
$t.Cell(12,1).Range.Text = "0.09883"

# Row 44 (0-based row 43): multi-run "103\t0.00019\t...\t100.0" -> single "99.89"
$t.Cell(44,1).Range.Text = "99.89"

# Row 45 (0-based row 44): multi-run "103\t0.00011\t...\t100.0" -> single "0.1"
$t.Cell(45,1).Range.Text = "0.1"

# Row 46 (0-based row 45): multi-run "103\t0.00011\t...\t100.0" -> single "90"
$t.Cell(46,1).Range.Text = "90"
